# Generate Report for Handoff
#
# The workbook tracks a "1e81a8ee-5a5a-415e-b176-6833c9bdf119" source file
# alongside "0c41f34b-3624-4887-a2e2-b9ddf5f65ac0". That source file's row
# is removed from every sheet (it's no longer part of the handoff set), the
# status of the remaining file flips from "Handed back: in sync with en-US"
# to "Ready for handoff", and the corresponding "Latest Handoff Datetime"
# stamps are refreshed to reflect the new handoff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# Row 2 = 0c41f34b...md, Row 3 = 1e81a8ee...md (to remove), Row 4 = .localization-config
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Clear hyperlinks first -- this engine doesn't re-target hyperlinks when
# rows shift, so they're rebuilt from scratch after the row is removed.
$wsOverview.Hyperlinks.Delete()

# Drop the 1e81a8ee row entirely; row 4 (.localization-config) shifts to row 3.
$wsOverview.Rows.Item(3).Delete()

# Status text for the remaining file is now "ready for handoff".
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b7e8acbb10cf77b36dad46156c87000dd916e617/e2e/0c41f34b-3624-4887-a2e2-b9ddf5f65ac0.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b7e8acbb10cf77b36dad46156c87000dd916e617/.localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": detail table for the zh-cn locale
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()
$wsZh.Rows.Item(3).Delete()

$wsZh.Range("B2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = "2016-03-09 15:34:16"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b7e8acbb10cf77b36dad46156c87000dd916e617/e2e/0c41f34b-3624-4887-a2e2-b9ddf5f65ac0.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc0121b0c25a4d77c91acf0a10308fcfd68d04cb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0c41f34b-3624-4887-a2e2-b9ddf5f65ac0.36ab0efcb4ee86b6f11b521db688b0665d165b93.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ef2bda5e02e42de8eb0b30f63a095df1a72aa31f/e2e/0c41f34b-3624-4887-a2e2-b9ddf5f65ac0.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/aaf966088b5892e09f9cbc875a1c2a44b22a01a4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0c41f34b-3624-4887-a2e2-b9ddf5f65ac0.36ab0efcb4ee86b6f11b521db688b0665d165b93.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b7e8acbb10cf77b36dad46156c87000dd916e617/.localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": detail table for the de-de locale
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()
$wsDe.Rows.Item(3).Delete()

$wsDe.Range("B2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = "2016-03-09 15:34:25"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b7e8acbb10cf77b36dad46156c87000dd916e617/e2e/0c41f34b-3624-4887-a2e2-b9ddf5f65ac0.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0af8dfaacf671c985302a153641efd37d190c2f3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0c41f34b-3624-4887-a2e2-b9ddf5f65ac0.36ab0efcb4ee86b6f11b521db688b0665d165b93.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1c00a2abbcc493de39efd83f2cad27d3b126d602/e2e/0c41f34b-3624-4887-a2e2-b9ddf5f65ac0.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/098d9ccab71dcd6bf3fc9e547b3d24a782b4e85a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0c41f34b-3624-4887-a2e2-b9ddf5f65ac0.36ab0efcb4ee86b6f11b521db688b0665d165b93.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b7e8acbb10cf77b36dad46156c87000dd916e617/.localization-config") | Out-Null
